$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$ws.Range("B3").Value = "3.8.0"

# Update Date value (row 8)
$ws.Range("B8").Value = "2022-08-12T09:44:57-05:00"

# Clear Copyright value (row 16)
$ws.Range("B16").Value = ""
